# The deck currently uses the "Integral" theme (ppt/theme/theme2.xml) for its
# slide master / slides, while an unused "Office Theme" colour set lives in
# ppt/theme/theme1.xml (only referenced by the notes master). The authored
# change swaps the content of the two theme parts, so the slides end up
# rendered with the default "Office Theme" colours instead of "Integral".
#
# The PowerPoint COM object model only exposes the *active* theme (the one
# backing the slide master / slides) for editing, via
# Master.ColorScheme.Colors(index).RGB — so we reproduce the visible part of
# that swap by repainting every theme colour slot on the slide master with
# the "Office Theme" palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.ColorScheme

# VBA RGB() packs bytes as 0xBBGGRR, so the values below are the reversed
# byte order of the target "Office Theme" hex colours.
$cs.Colors(1).RGB  = 0x000000   # dk1      -> 000000
$cs.Colors(2).RGB  = 0xFFFFFF   # lt1      -> FFFFFF
$cs.Colors(3).RGB  = 0x6A5444   # dk2      -> 44546A
$cs.Colors(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$cs.Colors(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$cs.Colors(6).RGB  = 0x317DED   # accent2  -> ED7D31
$cs.Colors(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$cs.Colors(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$cs.Colors(9).RGB  = 0xC47244   # accent5  -> 4472C4
$cs.Colors(10).RGB = 0x47AD70   # accent6  -> 70AD47
$cs.Colors(11).RGB = 0xC16305   # hlink    -> 0563C1
$cs.Colors(12).RGB = 0x724F95   # folHlink -> 954F72
